# The sheet previously held a single column (A1:A23) of numbers/words used
# for some lookup exercise. Replace it with a small word-list table:
#   column A = word, column B = a "1" marker (omitted on the last row),
#   column C = IF(B<>1,A,) formula mirroring A whenever B isn't 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @("cat", "dog", "mouse", "lock", "hourse", "mouse", "click", "damn", "damp", "clock")

# Wipe any leftover data from the old, longer single-column list so the
# used range shrinks back down to the new table.
$ws.Range("A1:C23").ClearContents()

for ($i = 0; $i -lt $words.Length; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $words[$i]
    if ($r -lt $words.Length) {
        $ws.Cells.Item($r, 2).Value = 1
    }
    $ws.Cells.Item($r, 3).Formula = "=IF(B$r<>1,A$r,)"
}

$ws.Range("A1").Select() | Out-Null
